# Daily attendance processing - 2025-10-12 10:43:21
# Normalizes the "Recorded By" (col G) author ordering so "System" is
# listed first, and applies two attendance-count corrections (H57, S17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Setting .Value directly on a cell whose new text looks like a number
# (e.g. "60.7%") makes Excel reinterpret it as a numeric percentage and
# reformat the cell. Routing the write through a quoted text formula and
# then collapsing it to a static value with Copy/PasteSpecial(xlPasteValues)
# keeps every target cell a literal text string in its original style.

$ws.Range("G2").Formula = "=""System, backup@backdoor.com, system"""
$ws.Range("G2").Copy()
$ws.Range("G2").PasteSpecial(-4163)
$ws.Range("G3").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G3").Copy()
$ws.Range("G3").PasteSpecial(-4163)
$ws.Range("G4").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G4").Copy()
$ws.Range("G4").PasteSpecial(-4163)
$ws.Range("G5").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G5").Copy()
$ws.Range("G5").PasteSpecial(-4163)
$ws.Range("G6").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G6").Copy()
$ws.Range("G6").PasteSpecial(-4163)
$ws.Range("G7").Formula = "=""System, admin@admin.com"""
$ws.Range("G7").Copy()
$ws.Range("G7").PasteSpecial(-4163)
$ws.Range("G10").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G10").Copy()
$ws.Range("G10").PasteSpecial(-4163)
$ws.Range("G11").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G11").Copy()
$ws.Range("G11").PasteSpecial(-4163)
$ws.Range("G12").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G12").Copy()
$ws.Range("G12").PasteSpecial(-4163)
$ws.Range("G13").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G13").Copy()
$ws.Range("G13").PasteSpecial(-4163)
$ws.Range("G14").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G14").Copy()
$ws.Range("G14").PasteSpecial(-4163)
$ws.Range("G15").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G15").Copy()
$ws.Range("G15").PasteSpecial(-4163)
$ws.Range("G17").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G17").Copy()
$ws.Range("G17").PasteSpecial(-4163)
$ws.Range("S17").Formula = "=""60.7%"""
$ws.Range("S17").Copy()
$ws.Range("S17").PasteSpecial(-4163)
$ws.Range("G18").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G18").Copy()
$ws.Range("G18").PasteSpecial(-4163)
$ws.Range("G29").Formula = "=""System, backup@backdoor.com, system"""
$ws.Range("G29").Copy()
$ws.Range("G29").PasteSpecial(-4163)
$ws.Range("G30").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G30").Copy()
$ws.Range("G30").PasteSpecial(-4163)
$ws.Range("G31").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G31").Copy()
$ws.Range("G31").PasteSpecial(-4163)
$ws.Range("G32").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G32").Copy()
$ws.Range("G32").PasteSpecial(-4163)
$ws.Range("G33").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G33").Copy()
$ws.Range("G33").PasteSpecial(-4163)
$ws.Range("G34").Formula = "=""System, admin@admin.com"""
$ws.Range("G34").Copy()
$ws.Range("G34").PasteSpecial(-4163)
$ws.Range("G37").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G37").Copy()
$ws.Range("G37").PasteSpecial(-4163)
$ws.Range("G38").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G38").Copy()
$ws.Range("G38").PasteSpecial(-4163)
$ws.Range("G39").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G39").Copy()
$ws.Range("G39").PasteSpecial(-4163)
$ws.Range("G40").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G40").Copy()
$ws.Range("G40").PasteSpecial(-4163)
$ws.Range("G41").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G41").Copy()
$ws.Range("G41").PasteSpecial(-4163)
$ws.Range("G42").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G42").Copy()
$ws.Range("G42").PasteSpecial(-4163)
$ws.Range("G44").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G44").Copy()
$ws.Range("G44").PasteSpecial(-4163)
$ws.Range("G45").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G45").Copy()
$ws.Range("G45").PasteSpecial(-4163)
$ws.Range("G56").Formula = "=""System, backup@backdoor.com, system"""
$ws.Range("G56").Copy()
$ws.Range("G56").PasteSpecial(-4163)
$ws.Range("G57").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G57").Copy()
$ws.Range("G57").PasteSpecial(-4163)
$ws.Range("H57").Formula = "=""29/55"""
$ws.Range("H57").Copy()
$ws.Range("H57").PasteSpecial(-4163)
$ws.Range("G58").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G58").Copy()
$ws.Range("G58").PasteSpecial(-4163)
$ws.Range("G59").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G59").Copy()
$ws.Range("G59").PasteSpecial(-4163)
$ws.Range("G60").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G60").Copy()
$ws.Range("G60").PasteSpecial(-4163)
$ws.Range("G61").Formula = "=""System, admin@admin.com"""
$ws.Range("G61").Copy()
$ws.Range("G61").PasteSpecial(-4163)
$ws.Range("G64").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G64").Copy()
$ws.Range("G64").PasteSpecial(-4163)
$ws.Range("G65").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G65").Copy()
$ws.Range("G65").PasteSpecial(-4163)
$ws.Range("G66").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G66").Copy()
$ws.Range("G66").PasteSpecial(-4163)
$ws.Range("G67").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G67").Copy()
$ws.Range("G67").PasteSpecial(-4163)
$ws.Range("G68").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G68").Copy()
$ws.Range("G68").PasteSpecial(-4163)
$ws.Range("G69").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G69").Copy()
$ws.Range("G69").PasteSpecial(-4163)
$ws.Range("G71").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G71").Copy()
$ws.Range("G71").PasteSpecial(-4163)
$ws.Range("G72").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G72").Copy()
$ws.Range("G72").PasteSpecial(-4163)
$ws.Range("G84").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G84").Copy()
$ws.Range("G84").PasteSpecial(-4163)
$ws.Range("G85").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G85").Copy()
$ws.Range("G85").PasteSpecial(-4163)
$ws.Range("G86").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G86").Copy()
$ws.Range("G86").PasteSpecial(-4163)
$ws.Range("G87").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G87").Copy()
$ws.Range("G87").PasteSpecial(-4163)
$ws.Range("G88").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G88").Copy()
$ws.Range("G88").PasteSpecial(-4163)
$ws.Range("G89").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G89").Copy()
$ws.Range("G89").PasteSpecial(-4163)
$ws.Range("G90").Formula = "=""dnasr281@gmail.com, admin@admin.com"""
$ws.Range("G90").Copy()
$ws.Range("G90").PasteSpecial(-4163)
$ws.Range("G93").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G93").Copy()
$ws.Range("G93").PasteSpecial(-4163)
$ws.Range("G95").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G95").Copy()
$ws.Range("G95").PasteSpecial(-4163)
$ws.Range("G96").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G96").Copy()
$ws.Range("G96").PasteSpecial(-4163)
$ws.Range("G110").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G110").Copy()
$ws.Range("G110").PasteSpecial(-4163)
$ws.Range("G111").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G111").Copy()
$ws.Range("G111").PasteSpecial(-4163)
$ws.Range("G112").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G112").Copy()
$ws.Range("G112").PasteSpecial(-4163)
$ws.Range("G113").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G113").Copy()
$ws.Range("G113").PasteSpecial(-4163)
$ws.Range("G114").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G114").Copy()
$ws.Range("G114").PasteSpecial(-4163)
$ws.Range("G115").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G115").Copy()
$ws.Range("G115").PasteSpecial(-4163)
$ws.Range("G116").Formula = "=""dnasr281@gmail.com, admin@admin.com"""
$ws.Range("G116").Copy()
$ws.Range("G116").PasteSpecial(-4163)
$ws.Range("G119").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G119").Copy()
$ws.Range("G119").PasteSpecial(-4163)
$ws.Range("G121").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G121").Copy()
$ws.Range("G121").PasteSpecial(-4163)
$ws.Range("G122").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G122").Copy()
$ws.Range("G122").PasteSpecial(-4163)
$ws.Range("G136").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G136").Copy()
$ws.Range("G136").PasteSpecial(-4163)
$ws.Range("G137").Formula = "=""System, backup@backdoor.com"""
$ws.Range("G137").Copy()
$ws.Range("G137").PasteSpecial(-4163)
$ws.Range("G138").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G138").Copy()
$ws.Range("G138").PasteSpecial(-4163)
$ws.Range("G139").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G139").Copy()
$ws.Range("G139").PasteSpecial(-4163)
$ws.Range("G140").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G140").Copy()
$ws.Range("G140").PasteSpecial(-4163)
$ws.Range("G141").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G141").Copy()
$ws.Range("G141").PasteSpecial(-4163)
$ws.Range("G142").Formula = "=""dnasr281@gmail.com, admin@admin.com"""
$ws.Range("G142").Copy()
$ws.Range("G142").PasteSpecial(-4163)
$ws.Range("G145").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G145").Copy()
$ws.Range("G145").PasteSpecial(-4163)
$ws.Range("G147").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G147").Copy()
$ws.Range("G147").PasteSpecial(-4163)
$ws.Range("G148").Formula = "=""System, dnasr281@gmail.com"""
$ws.Range("G148").Copy()
$ws.Range("G148").PasteSpecial(-4163)
$excel.CutCopyMode = $false
